# All Hat And No Cattle.xlsx - minor content update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (Satirist): ability text rewritten
$ws.Range("C14").Value = "If you are chosen at night by a minion, their ability malfunctions; the choice is redirected (without warning) at a new player of your choice."

# Row 16 (old Guard row): content removed -> row becomes blank in B:D (F/I/J untouched)
$ws.Range("B16:D16").ClearContents()

# Row 26: Bazooka replaced by Enchanter
$ws.Range("B26").Value = "Enchanter"
$ws.Range("C26").Value = "If a player targets the Enchanter, they become enchanted until targeted by another player. Enchanted players register as the Enchanter. "
$ws.Range("D26").Value = "Enchanted"

# Row 27 (Siren): ability text rewritten
$ws.Range("C27").Value = "If you are nominated by an opposing player, they become drunk, even if you are dead."

# Row 31 and 32: Hox / Lumpus swap places with updated ability text
$ws.Range("B31").Value = "Hox"
$ws.Range("C31").Value = "Each night*, choose a player to die. Once per game, instead, choose to kill all players who were on the block today"

$ws.Range("B32").Value = "Lumpus"
$ws.Range("C32").Value = "Each night*, choose a player to die, even if they might have survived for some reason. After your first kill, your target must neighbour a dead player."
$ws.Range("D32").Value = "Dead, On the block"

# Row 35: now holds Guard (new ability text)
$ws.Range("B35").Value = "Guard"
$ws.Range("C35").Value = "Other players who voted on your exile today cannot die tonight (if you are not exiled), and one is drunk tonight."
$ws.Range("D35").Value = "Cannot die, Drunk"

# Row 36: now holds Jailer
$ws.Range("B36").Value = "Jailer"
$ws.Range("C36").Value = "Each night*, choose a player, they are placed on the block at dawn."

# Row 37: new row holding Triffid (moved down from row 36, updated text)
$ws.Range("B37").Value = "Triffid"
$ws.Range("C37").Value = "On a tied nomination vote, two (or more) players are placed on the block. At the end of the day, you choose which one of them is executed."
$ws.Range("D37").Value = "On the block"
$ws.Range("F37").Value = 0

# Update selection to match the saved view
$ws.Range("B18").Select()
